$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.959.44"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.219.38"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "291.85"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.28"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.515"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.466"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.47"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  +5.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0779"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.564.40"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.79"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.266.66"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.893.51"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.13"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.46"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.00"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.22"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -7.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.34"
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.91"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0994"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  -4.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.088.29"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0271"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.07"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.82"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  -8.71%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.434.19"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +2.40%  "
